$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'312.79"
$ws.Range("E2").Value = "'2.38%"
$ws.Range("E3").Value = "'1.34%"
$ws.Range("D4").Value = "'5.154"
$ws.Range("E4").Value = "'1.23%"
$ws.Range("D5").Value = "'0.07889"
$ws.Range("E5").Value = "'2.29%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.916"
$ws.Range("E6").Value = "'2.63%"
$ws.Range("E7").Value = "'1.10%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.853"
$ws.Range("E8").Value = "'-10.66%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9215"
$ws.Range("E9").Value = "'0.34%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1215"
$ws.Range("E10").Value = "'1.03%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1918"
$ws.Range("E11").Value = "'1.92%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09162"
$ws.Range("E12").Value = "'4.96%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03348"
$ws.Range("E13").Value = "'-1.07%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09586"
$ws.Range("E14").Value = "'-1.26%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001383"
$ws.Range("E15").Value = "'1.19%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005687"
$ws.Range("E16").Value = "'-6.92%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.511"
$ws.Range("E17").Value = "'-1.48%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.412"
$ws.Range("E18").Value = "'1.16%"
$ws.Range("D19").Value = "'0.3442"
$ws.Range("E19").Value = "'2.02%"
$ws.Range("D20").Value = "'5.258"
$ws.Range("E20").Value = "'4.56%"
$ws.Range("E21").Value = "'-0.26%"
$ws.Range("E22").Value = "'3.85%"
$ws.Range("E23").Value = "'-0.34%"
$ws.Range("D24").Value = "'0.04360"
$ws.Range("E24").Value = "'0.58%"
$ws.Range("D25").Value = "'0.001248"
$ws.Range("E25").Value = "'2.62%"
$ws.Range("D26").Value = "'0.004709"
$ws.Range("E26").Value = "'5.47%"
$ws.Range("E27").Value = "'-9.85%"
$ws.Range("D39").Value = "'0.02299"
$ws.Range("E39").Value = "'3.40%"
$ws.Range("D40").Value = "'0.05085"
$ws.Range("E40").Value = "'3.50%"
$ws.Range("D41").Value = "'0.007479"
$ws.Range("E41").Value = "'-1.59%"
$ws.Range("D42").Value = "'0.008883"
$ws.Range("E42").Value = "'-10.24%"
$ws.Range("D43").Value = "'0.1356"
$ws.Range("E43").Value = "'2.09%"
$ws.Range("D44").Value = "'0.001951"
$ws.Range("E44").Value = "'-5.74%"
$ws.Range("D45").Value = "'0.008623"
$ws.Range("E45").Value = "'-2.54%"
$ws.Range("D46").Value = "'0.00006618"
$ws.Range("E46").Value = "'-2.55%"
$ws.Range("E47").Value = "'-0.27%"
$ws.Range("D48").Value = "'0.003354"
$ws.Range("E48").Value = "'11.48%"
$ws.Range("E49").Value = "'-7.95%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.27%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.27%"
